# The commit updates the sample/template config table so that the "bu"
# field-name example reads "bu666" (the rest of the template data such as
# the demo type rows and the JSON sample value in C4 are left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("s1")
$ws.Range("D1").Value = "bu666"
